$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '26.799.24'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '1.638.22'
$ws.Range("E3").Value = '  +0.05%  '
$ws.Range("E4").Value = '  -0.61%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.90'
$ws.Range("E5").Value = '  -0.97%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.514'
$ws.Range("E6").Value = '  +3.10%  '
$ws.Range("E7").Value = '  -0.54%  '
$ws.Range("E8").Value = '  +1.80%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0623'
$ws.Range("E9").Value = '  +0.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.87'
$ws.Range("E10").Value = '  +3.55%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0844'
$ws.Range("E11").Value = '  -0.01%  '
$ws.Range("D12").Value = '1.866.70'
$ws.Range("E12").Value = '  -0.05%  '
$ws.Range("D13").Value = '1.645.09'
$ws.Range("E13").Value = '  +0.38%  '
$ws.Range("E14").Value = '  -0.29%  '
$ws.Range("E15").Value = '  +0.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.56'
$ws.Range("E16").Value = '  +3.08%  '
$ws.Range("D17").Value = '26.799.41'
$ws.Range("E17").Value = '  +0.09%  '
$ws.Range("D18").Value = '0.0₃0729'
$ws.Range("E18").Value = '  -0.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '219.06'
$ws.Range("E19").Value = '  +1.93%  '
$ws.Range("E20").Value = '  -0.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.71'
$ws.Range("E21").Value = '  +6.78%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.39'
$ws.Range("E23").Value = '  +3.64%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.15'
$ws.Range("E24").Value = '  +0.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.18'
$ws.Range("E25").Value = '  -0.39%  '
$ws.Range("E26").Value = '  -0.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.36'
$ws.Range("E27").Value = '  +4.49%  '
$ws.Range("E28").Value = '  +0.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.77'
$ws.Range("E29").Value = '  +0.62%  '
$ws.Range("E30").Value = '  -0.53%  '
$ws.Range("E31").Value = '  -1.74%  '
$ws.Range("E32").Value = '  -1.68%  '
$ws.Range("E33").Value = '  +0.61%  '
$ws.Range("E34").Value = '  +1.50%  '
$ws.Range("D35").Value = '1.258.32'
$ws.Range("E35").Value = '  -0.16%  '
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("E37").Value = '  +1.09%  '
$ws.Range("E38").Value = '  +0.99%  '
$ws.Range("E39").Value = '  +2.79%  '
$ws.Range("E40").Value = '  -0.45%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.807'
$ws.Range("E41").Value = '  +0.27%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.46'
$ws.Range("E42").Value = '  +2.80%  '
$ws.Range("D43").Value = '1.777.43'
$ws.Range("E43").Value = '  +0.10%  '
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '61.79'
$ws.Range("E44").Value = '  +3.00%  '
$ws.Range("B45").Value = 'MXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.10'
$ws.Range("E45").Value = '  -1.35%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.75'
$ws.Range("E46").Value = '  -0.36%  '
$ws.Range("E47").Value = '  -1.11%  '
$ws.Range("E48").Value = '  -0.63%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.61'
$ws.Range("E49").Value = '  +0.84%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0965'
$ws.Range("E50").Value = '  +0.49%  '
$ws.Range("E51").Value = '  -0.52%  '
